$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where D/J/K/L/M/P values need updating (values are moved between rows).
# Mapping built from diff: for each destination row, the new tuple
# (Fecha, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado, Precio $/Kg)

$updates = @{
    2  = @(44284, 1500, 800, 850, 825, 825)
    4  = @(44607, 900, 1300, 1400, 1350, 1350)
    5  = @(44638, 1000, 900, 950, 925, 925)
    6  = @(44229, 1500, 1400, 1500, 1450, 1450)
    7  = @(44291, 1000, 1000, 1200, 1100, 1100)
    8  = @(44656, 1000, 900, 1000, 950, 950)
    9  = @(44175, 1600, 1000, 1200, 1100, 1100)
    11 = @(44673, 900, 1300, 1400, 1350, 1350)
    12 = @(44341, 1300, 900, 1000, 950, 950)
    13 = @(44455, 1100, 900, 1000, 950, 950)
    14 = @(44243, 1200, 1200, 1300, 1250, 1250)
    15 = @(44453, 1000, 800, 900, 850, 850)
    16 = @(44407, 1000, 1200, 1300, 1250, 1250)
    17 = @(44476, 900, 700, 800, 750, 750)
    18 = @(44649, 600, 900, 1000, 950, 950)
    19 = @(44550, 1300, 1000, 1200, 1100, 1100)
    20 = @(44449, 1300, 900, 950, 925, 925)
    21 = @(44687, 1000, 1200, 1300, 1250, 1250)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
    $ws.Cells.Item($row, 11).Value = $vals[2]
    $ws.Cells.Item($row, 12).Value = $vals[3]
    $ws.Cells.Item($row, 13).Value = $vals[4]
    $ws.Cells.Item($row, 16).Value = $vals[5]
}
